$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 143, shifting existing rows 143-159 down to 144-160.
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with the new weekly data record.
$ws.Range("A143").Value = 6
$ws.Range("B143").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C143").Value = "Metropolitana"
$ws.Range("D143").Value = 44491
$ws.Range("D143").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = 100112026
$ws.Range("G143").Value = "Haba"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 800
$ws.Range("K143").Value = 6000
$ws.Range("L143").Value = 7000
$ws.Range("M143").Value = 6438
$ws.Range("N143").Value = '$/saco 25 kilos'
$ws.Range("O143").Value = "Provincia de Melipilla"
$ws.Range("P143").Value = 258
$ws.Range("Q143").Value = 25
$ws.Range("R143").Value = "Hortaliza"
